# 21-22 KBL stats: replace the short team abbreviation in column B with the
# full "city + team" name for every player row (rows 2-184), grouped by the
# contiguous block of rows that belong to each team.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2:B18").Value   = "서울 SK"
$ws.Range("B19:B37").Value  = "수원 KT"
$ws.Range("B38:B56").Value  = "부산 KCC"
$ws.Range("B57:B74").Value  = "울산 모비스"
$ws.Range("B75:B95").Value  = "원주 DB"
$ws.Range("B96:B114").Value = "고양 소노"
$ws.Range("B115:B132").Value = "서울 삼성"
$ws.Range("B133:B150").Value = "창원 LG"
$ws.Range("B151:B166").Value = "대구 가스공사"
$ws.Range("B167:B184").Value = "안양 KGC"
